$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.045.52'
$ws.Range("E2").Value = '  -3.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.364.59'
$ws.Range("E3").Value = '  -3.94%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '502.04'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.70'
$ws.Range("E6").Value = '  -3.09%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  -2.74%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.368.16'
$ws.Range("E9").Value = '  -3.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0983'
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.80'
$ws.Range("E12").Value = '  +3.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.785.14'
$ws.Range("E14").Value = '  -3.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.034.07'
$ws.Range("E15").Value = '  -3.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.44'
$ws.Range("E16").Value = '  -2.44%  '
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.339.45'
$ws.Range("E18").Value = '  -5.47%  '
$ws.Range("E19").Value = '  -3.43%  '
$ws.Range("E20").Value = '  -2.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '307.61'
$ws.Range("E21").Value = '  -2.31%  '
$ws.Range("E22").Value = '  -2.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.88'
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.997'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("E26").Value = '  -3.10%  '
$ws.Range("E27").Value = '  -6.09%  '
$ws.Range("E28").Value = '  -5.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.41'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0710'
$ws.Range("E30").Value = '  -3.57%  '
$ws.Range("E31").Value = '  -2.95%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.79'
$ws.Range("E33").Value = '  -6.11%  '
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("E35").Value = '  -5.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.60'
$ws.Range("E36").Value = '  -2.77%  '
$ws.Range("E37").Value = '  -5.96%  '
$ws.Range("E38").Value = '  -4.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.18'
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.801'
$ws.Range("E40").Value = '  -0.99%  '
$ws.Range("E41").Value = '  -5.94%  '
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '129.13'
$ws.Range("E43").Value = '  -5.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.70'
$ws.Range("E44").Value = '  -4.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.562'
$ws.Range("E45").Value = '  -2.87%  '
$ws.Range("E46").Value = '  -1.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '238.68'
$ws.Range("E47").Value = '  -6.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0482'
$ws.Range("E48").Value = '  -2.50%  '
$ws.Range("E49").Value = '  -3.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.98'
$ws.Range("E50").Value = '  -2.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.949'
$ws.Range("E51").Value = '  -1.25%  '
